# "Generate Report for Archive" - refresh the localization-status report:
# the "Ready for handoff" status has moved on to "In Translation" for both
# tracked files, on the Overview sheet as well as the per-language sheets.
# Excel auto-shrinks the Status columns to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Overview sheet: status is mirrored per-language in columns E (zh-cn) and F (de-de)
foreach ($cellRef in @("E2", "F2", "E3", "F3")) {
    $cell = $wsOverview.Range($cellRef)
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# Per-language sheets: status lives in column C
foreach ($ws in @($wsZhCn, $wsDeDe)) {
    foreach ($cellRef in @("C2", "C3")) {
        $cell = $ws.Range($cellRef)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# The Status columns auto-fit to the new, narrower text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.576851254417766
$wsOverview.Columns.Item(6).ColumnWidth = 12.576851254417766
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.576851254417766
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.576851254417766
